$p = $ppt.ActivePresentation
$p.Slides.Item(6).Delete()
